$d = $word.ActiveDocument
$p1 = $d.Paragraphs(1)

# Add a (line-less) paragraph border with 5pt spacing on all four sides.
$borders = $p1.Range.ParagraphFormat.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5

# Widen the left indent from 120 twips (6pt) to 225 twips (11.25pt).
$p1.Range.ParagraphFormat.LeftIndent = 11.25

# Update the placeholder ID text; the identically-formatted trailing run merges
# into this one automatically once the text matches.
$find = $p1.Range.Find
$find.Execute("**ID__AFFARS_pgi_5306_topic_12__ID**", $true, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_USAFA_PGI_5306__ID**", 2) | Out-Null

# Remove the now-trailing single space left over after the replacement, so the
# paragraph holds only the placeholder text (no second run).
$newIdLen = "**ID__AFFARS_USAFA_PGI_5306__ID**".Length
$trailing = $d.Range($p1.Range.Start + $newIdLen, $p1.Range.Start + $newIdLen + 1)
if ($trailing.Text -eq " ") {
    $trailing.Delete() | Out-Null
}
